$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  "Tử Vi đồng cung với Tham Lang tại Phụ Mẫu",
  "Tử Vi đồng cung Thiên Phủ tại Phụ Mẫu",
  "Tử Vi đồng cung với Thiên Tướng tại Phụ Mẫu",
  "Tử Vi đồng cung với Thất Sát tại Phụ Mẫu",
  "Tử Vi đồng cung với Phá Quân tại Phụ Mẫu",
  "Liêm Trinh đồng cung với Thiên Phủ tại Phụ Mẫu",
  "Liêm Trinh đồng cung với Thiên Tướng tại Phụ Mẫu",
  "Liêm Trinh đồng cung với Thất Sát tại Phụ Mẫu",
  "Liêm Trinh đồng cung với Tham Lang tại Phụ Mẫu",
  "Thiên Đồng đồng cung với Thiên Lương tại Phụ Mẫu",
  "Thiên Đồng đồng cung với Thái Âm tại Phụ Mẫu ở Tý",
  "Thiên Đồng đồng cung với Thái Âm tại Phụ Mẫu ở Ngọ",
  "Thiên Đồng đồng cung với Cự Môn tại Phụ Mẫu",
  "Vũ Khúc tọa thủ cung Phụ Mẫu đồng cung Phá Quân",
  "Vũ Khúc tọa thủ cung Phụ Mẫu đồng cung Thất Sát",
  "Vũ Khúc tọa thủ cung Phụ Mẫu đồng cung Thiên Tướng",
  "Vũ Khúc tọa thủ cung Phụ Mẫu đồng cung Thiên Phủ",
  "Thái Dương tọa thủ cung Phụ Mẫu ở Sửu đồng cung Thái Âm",
  "Thái Dương tọa thủ cung Phụ Mẫu ở Mùi đồng cung Thái Âm",
  "Thiên Cơ tọa thủ cung Phụ Mẫu đồng cung Thái Âm ở Dần",
  "Thiên Cơ tọa thủ cung Phụ Mẫu đồng cung Thái Âm ở Thân",
  "Thiên Cơ tọa thủ cung Phụ Mẫu đồng cung Thiên Lương",
  "Thiên Cơ tọa thủ cung Phụ Mẫu đồng cung Cự Môn",
  "Tham Lang Vũ Khúc đồng cung tại cung Phụ Mẫu",
  "Thiên Lương đồng cung Thái Dương tại cung Phụ Mẫu ở Mão",
  "Thiên Lương đồng cung Thái Dương tại cung Phụ Mẫu ở Dậu",
  "Kình Dương, Liêm Trinh đồng cung tại Phụ Mẫu",
  "Kình Dương, Thất Sát đồng cung tại Phụ Mẫu",
  "Kình Dương, Tham Lang đồng cung tại Phụ Mẫu",
  "Đà La, Liêm Trinh đồng cung tại Phụ Mẫu",
  "Đà La, Thất Sát đồng cung tại Phụ Mẫu",
  "Đà La, Tham Lang đồng cung tại Phụ Mẫu",
  "Tham Lang đồng cung với Hỏa Tinh tại Phụ Mẫu",
  "Tham Lang đồng cung với Linh Tinh tại Phụ Mẫu",
  "Phá Quân đồng cung với Hỏa Tinh tại Phụ Mẫu",
  "Phá Quân đồng cung với Linh Tinh tại Phụ Mẫu",
  "Thái Dương, Thái Âm, Hóa kỵ đồng cung tại Sửu",
  "Thái Dương, Thái Âm, Hóa kỵ đồng cung tại Mùi",
  "Kình Dương, Đà La, Thiên Mã đồng cung tại Phụ Mẫu gặp Hỏa Tinh, Linh Tinh",
  "Kình Dương, Đà La, Thái Tuế đồng cung tại Phụ Mẫu gặp Hỏa Tinh, Linh Tinh",
  "Đào Hoa, Hồng Loan, Thái Tuế đồng cung tại Phụ Mẫu",
)

$startRow = 3815
for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 1).Value = $values[$i]
  $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("R3827").Select()

